$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1833.7333  # ALC!H19
$ws.Cells.Item(19, 10).Value = 1421.5  # ALC!J19
$ws.Cells.Item(19, 12).Value = 1421.5  # ALC!L19
$ws.Cells.Item(19, 14).Value = -1771.5  # ALC!N19

$ws.Cells.Item(33, 8).Value = 267  # ALC!H33
$ws.Cells.Item(33, 9).Value = 297.5  # ALC!I33
$ws.Cells.Item(33, 10).Value = 96.2  # ALC!J33
$ws.Cells.Item(33, 11).Value = 297.5  # ALC!K33
$ws.Cells.Item(33, 12).Value = 96.2  # ALC!L33
$ws.Cells.Item(33, 13).Value = -68.5  # ALC!M33
$ws.Cells.Item(33, 14).Value = -554.2  # ALC!N33

$ws.Cells.Item(40, 8).Value = 1536.0476  # ALC!H40
$ws.Cells.Item(40, 9).Value = 1139.8  # ALC!I40
$ws.Cells.Item(40, 10).Value = 1659.875  # ALC!J40
$ws.Cells.Item(40, 11).Value = 1139.8  # ALC!K40
$ws.Cells.Item(40, 12).Value = 1659.875  # ALC!L40
$ws.Cells.Item(40, 13).Value = -964.8  # ALC!M40
$ws.Cells.Item(40, 14).Value = -2009.875  # ALC!N40

$ws.Cells.Item(125, 8).Value = 6000  # ALC!H125
$ws.Cells.Item(125, 9).Value = 10000  # ALC!I125
$ws.Cells.Item(125, 10).Value = 2000  # ALC!J125
$ws.Cells.Item(125, 11).Value = 90000  # ALC!K125
$ws.Cells.Item(125, 12).Value = 18000  # ALC!L125
$ws.Cells.Item(125, 13).Value = -87540  # ALC!M125
$ws.Cells.Item(125, 14).Value = -22920  # ALC!N125

$ws.Cells.Item(140, 8).Value = 48716.668  # ALC!H140
$ws.Cells.Item(140, 10).Value = 48716.668  # ALC!J140
$ws.Cells.Item(140, 12).Value = 48716.668  # ALC!L140
$ws.Cells.Item(140, 14).Value = -59076.668  # ALC!N140

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 279  # ARM!H26
$ws.Cells.Item(26, 9).Value = 279  # ARM!I26
$ws.Cells.Item(26, 11).Value = 279  # ARM!K26
$ws.Cells.Item(26, 13).Value = 51  # ARM!M26

$ws.Cells.Item(28, 8).Value = 3614  # ARM!H28
$ws.Cells.Item(28, 9).Value = 3614  # ARM!I28
$ws.Cells.Item(28, 11).Value = 3614  # ARM!K28
$ws.Cells.Item(28, 13).Value = -3422  # ARM!M28

$ws.Cells.Item(61, 8).Value = 55560000  # ARM!H61
$ws.Cells.Item(61, 9).Value = 62504628  # ARM!I61
$ws.Cells.Item(61, 10).Value = 3000  # ARM!J61
$ws.Cells.Item(61, 11).Value = 62504628  # ARM!K61
$ws.Cells.Item(61, 12).Value = 3000  # ARM!L61
$ws.Cells.Item(61, 13).Value = -62504416  # ARM!M61
$ws.Cells.Item(61, 14).Value = -3424  # ARM!N61

$ws.Cells.Item(99, 8).Value = 3614  # ARM!H99
$ws.Cells.Item(99, 9).Value = 3614  # ARM!I99
$ws.Cells.Item(99, 11).Value = 3614  # ARM!K99
$ws.Cells.Item(99, 13).Value = -619  # ARM!M99

$ws.Cells.Item(132, 8).Value = 11908272  # ARM!H132
$ws.Cells.Item(132, 9).Value = 22730254  # ARM!I132
$ws.Cells.Item(132, 10).Value = 4092.6  # ARM!J132
$ws.Cells.Item(132, 11).Value = 68190762  # ARM!K132
$ws.Cells.Item(132, 12).Value = 12277.8  # ARM!L132
$ws.Cells.Item(132, 13).Value = -68188232  # ARM!M132
$ws.Cells.Item(132, 14).Value = -17337.8  # ARM!N132

$ws.Cells.Item(136, 8).Value = 55560000  # ARM!H136
$ws.Cells.Item(136, 9).Value = 62504628  # ARM!I136
$ws.Cells.Item(136, 10).Value = 3000  # ARM!J136
$ws.Cells.Item(136, 11).Value = 187513884  # ARM!K136
$ws.Cells.Item(136, 12).Value = 9000  # ARM!L136
$ws.Cells.Item(136, 13).Value = -187511334  # ARM!M136
$ws.Cells.Item(136, 14).Value = -14100  # ARM!N136

$ws.Cells.Item(141, 8).Value = 40666.668  # ARM!H141
$ws.Cells.Item(141, 10).Value = 40666.668  # ARM!J141
$ws.Cells.Item(141, 12).Value = 40666.668  # ARM!L141
$ws.Cells.Item(141, 14).Value = -51026.668  # ARM!N141

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 1003  # BSM!H7
$ws.Cells.Item(7, 9).Value = 1003  # BSM!I7
$ws.Cells.Item(7, 10).Value = 0  # BSM!J7
$ws.Cells.Item(7, 11).Value = 1003  # BSM!K7
$ws.Cells.Item(7, 12).Value = 0  # BSM!L7
$ws.Cells.Item(7, 13).Value = -890  # BSM!M7
$ws.Cells.Item(7, 14).ClearContents()  # BSM!N7

$ws.Cells.Item(11, 8).Value = 1004.0833  # BSM!H11
$ws.Cells.Item(11, 9).Value = 1010.5714  # BSM!I11
$ws.Cells.Item(11, 10).Value = 995  # BSM!J11
$ws.Cells.Item(11, 11).Value = 1010.5714  # BSM!K11
$ws.Cells.Item(11, 12).Value = 995  # BSM!L11
$ws.Cells.Item(11, 13).Value = -870.5714  # BSM!M11
$ws.Cells.Item(11, 14).Value = -1275  # BSM!N11

$ws.Cells.Item(134, 8).Value = 3962.3142  # BSM!H134
$ws.Cells.Item(134, 9).Value = 3152.4614  # BSM!I134
$ws.Cells.Item(134, 10).Value = 6301.8887  # BSM!J134
$ws.Cells.Item(134, 11).Value = 9457.3842  # BSM!K134
$ws.Cells.Item(134, 12).Value = 18905.6661  # BSM!L134
$ws.Cells.Item(134, 13).Value = -6922.3842  # BSM!M134
$ws.Cells.Item(134, 14).Value = -23975.6661  # BSM!N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7249870.5  # CRP!H31
$ws.Cells.Item(31, 9).Value = 3599.524  # CRP!I31
$ws.Cells.Item(31, 10).Value = 83335720  # CRP!J31
$ws.Cells.Item(31, 11).Value = 3599.524  # CRP!K31
$ws.Cells.Item(31, 12).Value = 83335720  # CRP!L31
$ws.Cells.Item(31, 13).Value = -3304.524  # CRP!M31
$ws.Cells.Item(31, 14).Value = -83336310  # CRP!N31

$ws.Cells.Item(34, 8).Value = 7249870.5  # CRP!H34
$ws.Cells.Item(34, 9).Value = 3599.524  # CRP!I34
$ws.Cells.Item(34, 10).Value = 83335720  # CRP!J34
$ws.Cells.Item(34, 11).Value = 3599.524  # CRP!K34
$ws.Cells.Item(34, 12).Value = 83335720  # CRP!L34
$ws.Cells.Item(34, 13).Value = -3397.524  # CRP!M34
$ws.Cells.Item(34, 14).Value = -83336124  # CRP!N34

$ws.Cells.Item(63, 8).Value = 36635.5  # CRP!H63
$ws.Cells.Item(63, 10).Value = 36635.5  # CRP!J63
$ws.Cells.Item(63, 12).Value = 36635.5  # CRP!L63
$ws.Cells.Item(63, 14).Value = -38007.5  # CRP!N63

$ws.Cells.Item(66, 8).Value = 36635.5  # CRP!H66
$ws.Cells.Item(66, 10).Value = 36635.5  # CRP!J66
$ws.Cells.Item(66, 12).Value = 109906.5  # CRP!L66
$ws.Cells.Item(66, 14).Value = -116770.5  # CRP!N66

$ws.Cells.Item(99, 8).Value = 1516.5  # CRP!H99
$ws.Cells.Item(99, 9).Value = 1339.6  # CRP!I99
$ws.Cells.Item(99, 10).Value = 1642.8572  # CRP!J99
$ws.Cells.Item(99, 11).Value = 1339.6  # CRP!K99
$ws.Cells.Item(99, 12).Value = 1642.8572  # CRP!L99
$ws.Cells.Item(99, 13).Value = 158.4000000000001  # CRP!M99
$ws.Cells.Item(99, 14).Value = -4638.8572  # CRP!N99

$ws.Cells.Item(126, 8).Value = 1516.5  # CRP!H126
$ws.Cells.Item(126, 9).Value = 1339.6  # CRP!I126
$ws.Cells.Item(126, 10).Value = 1642.8572  # CRP!J126
$ws.Cells.Item(126, 11).Value = 4018.8  # CRP!K126
$ws.Cells.Item(126, 12).Value = 4928.571599999999  # CRP!L126
$ws.Cells.Item(126, 13).Value = -1548.8  # CRP!M126
$ws.Cells.Item(126, 14).Value = -9868.571599999999  # CRP!N126

$ws.Cells.Item(134, 8).Value = 2367.6875  # CRP!H134
$ws.Cells.Item(134, 9).Value = 2569.9092  # CRP!I134
$ws.Cells.Item(134, 10).Value = 1922.8  # CRP!J134
$ws.Cells.Item(134, 11).Value = 7709.7276  # CRP!K134
$ws.Cells.Item(134, 12).Value = 5768.4  # CRP!L134
$ws.Cells.Item(134, 13).Value = -5174.7276  # CRP!M134
$ws.Cells.Item(134, 14).Value = -10838.4  # CRP!N134

$ws.Cells.Item(140, 8).Value = 42406.152  # CRP!H140
$ws.Cells.Item(140, 10).Value = 42406.152  # CRP!J140
$ws.Cells.Item(140, 12).Value = 42406.152  # CRP!L140
$ws.Cells.Item(140, 14).Value = -52766.152  # CRP!N140

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 33.904762  # CUL!H12
$ws.Cells.Item(12, 9).Value = 14.285714  # CUL!I12
$ws.Cells.Item(12, 10).Value = 43.714287  # CUL!J12
$ws.Cells.Item(12, 11).Value = 42.857142  # CUL!K12
$ws.Cells.Item(12, 12).Value = 131.142861  # CUL!L12
$ws.Cells.Item(12, 13).Value = 130.142858  # CUL!M12
$ws.Cells.Item(12, 14).Value = -477.142861  # CUL!N12

$ws.Cells.Item(70, 8).Value = 4642  # CUL!H70
$ws.Cells.Item(70, 10).Value = 5102.222  # CUL!J70
$ws.Cells.Item(70, 12).Value = 15306.666  # CUL!L70
$ws.Cells.Item(70, 14).Value = -15936.666  # CUL!N70

$ws.Cells.Item(73, 8).Value = 4642  # CUL!H73
$ws.Cells.Item(73, 10).Value = 5102.222  # CUL!J73
$ws.Cells.Item(73, 12).Value = 15306.666  # CUL!L73
$ws.Cells.Item(73, 14).Value = -17490.666  # CUL!N73

$ws.Cells.Item(95, 8).Value = 8660  # CUL!H95
$ws.Cells.Item(95, 9).Value = 9000  # CUL!I95
$ws.Cells.Item(95, 10).Value = 8490  # CUL!J95
$ws.Cells.Item(95, 11).Value = 27000  # CUL!K95
$ws.Cells.Item(95, 12).Value = 25470  # CUL!L95
$ws.Cells.Item(95, 13).Value = -24941  # CUL!M95
$ws.Cells.Item(95, 14).Value = -29588  # CUL!N95

$ws.Cells.Item(132, 8).Value = 821  # CUL!H132
$ws.Cells.Item(132, 9).Value = 529.8333  # CUL!I132
$ws.Cells.Item(132, 10).Value = 1519.8  # CUL!J132
$ws.Cells.Item(132, 11).Value = 4768.4997  # CUL!K132
$ws.Cells.Item(132, 12).Value = 13678.2  # CUL!L132
$ws.Cells.Item(132, 13).Value = -2238.4997  # CUL!M132
$ws.Cells.Item(132, 14).Value = -18738.2  # CUL!N132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 1250252.2  # GSM!H3
$ws.Cells.Item(3, 9).Value = 1250252.2  # GSM!I3
$ws.Cells.Item(3, 11).Value = 1250252.2  # GSM!K3
$ws.Cells.Item(3, 13).Value = -1250136.2  # GSM!M3

$ws.Cells.Item(102, 8).Value = 3110.577  # GSM!H102
$ws.Cells.Item(102, 9).Value = 3613.8  # GSM!I102
$ws.Cells.Item(102, 11).Value = 3613.8  # GSM!K102
$ws.Cells.Item(102, 13).Value = -1991.8  # GSM!M102

$ws.Cells.Item(122, 8).Value = 4168932  # GSM!H122
$ws.Cells.Item(122, 9).Value = 9525412  # GSM!I122
$ws.Cells.Item(122, 10).Value = 2780.6667  # GSM!J122
$ws.Cells.Item(122, 11).Value = 28576236  # GSM!K122
$ws.Cells.Item(122, 12).Value = 8342.000100000001  # GSM!L122
$ws.Cells.Item(122, 13).Value = -28573786  # GSM!M122
$ws.Cells.Item(122, 14).Value = -13242.0001  # GSM!N122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 166.66667  # LTW!H25
$ws.Cells.Item(25, 9).Value = 166.66667  # LTW!I25
$ws.Cells.Item(25, 11).Value = 166.66667  # LTW!K25
$ws.Cells.Item(25, 13).Value = 63.33332999999999  # LTW!M25

$ws.Cells.Item(56, 8).Value = 18064.428  # LTW!H56
$ws.Cells.Item(56, 9).Value = 13570.2  # LTW!I56
$ws.Cells.Item(56, 11).Value = 13570.2  # LTW!K56
$ws.Cells.Item(56, 13).Value = -12879.2  # LTW!M56

$ws.Cells.Item(132, 8).Value = 8797.306  # LTW!H132
$ws.Cells.Item(132, 9).Value = 5716.5625  # LTW!I132
$ws.Cells.Item(132, 10).Value = 11261.9  # LTW!J132
$ws.Cells.Item(132, 11).Value = 17149.6875  # LTW!K132
$ws.Cells.Item(132, 12).Value = 33785.7  # LTW!L132
$ws.Cells.Item(132, 13).Value = -14619.6875  # LTW!M132
$ws.Cells.Item(132, 14).Value = -38845.7  # LTW!N132

$ws.Cells.Item(136, 8).Value = 13894428  # LTW!H136
$ws.Cells.Item(136, 9).Value = 20002312  # LTW!I136
$ws.Cells.Item(136, 10).Value = 12873.182  # LTW!J136
$ws.Cells.Item(136, 11).Value = 60006936  # LTW!K136
$ws.Cells.Item(136, 12).Value = 38619.546  # LTW!L136
$ws.Cells.Item(136, 13).Value = -60004386  # LTW!M136
$ws.Cells.Item(136, 14).Value = -43719.546  # LTW!N136

$ws.Cells.Item(139, 8).Value = 59963.43  # LTW!H139
$ws.Cells.Item(139, 10).Value = 59963.43  # LTW!J139
$ws.Cells.Item(139, 12).Value = 59963.43  # LTW!L139
$ws.Cells.Item(139, 14).Value = -70243.42999999999  # LTW!N139
